$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for column F (dSF) per the repull / push all data / mean calculation.
$updates = @{
    2  = -1
    3  = -4
    4  = 5
    5  = -3
    6  = -5
    7  = -7
    8  = -1
    13 = -4
    14 = 2
    16 = -4
    17 = -4
    19 = -6
    20 = -5
    22 = 2
    24 = 2
    25 = 1
    26 = -5
    27 = -2
    28 = 2
    29 = -6
    30 = -3
    31 = -5
    32 = -1
    33 = 6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
